# Applies the odds updates described in the commit diff
# ("Atualizando o arquivo XLSX") to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("AE2").Value = 21
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 34
$ws.Range("AR2").Value = 67
$ws.Range("AU2").Value = 9.5
$ws.Range("BB2").Value = 451

# Row 3
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62

# Row 7
$ws.Range("G7").Value = 1.73
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.4
$ws.Range("L7").Value = 6
$ws.Range("Y7").Value = 9
$ws.Range("AA7").Value = 17
$ws.Range("AD7").Value = 7
$ws.Range("AJ7").Value = 19
$ws.Range("AL7").Value = 51
$ws.Range("AO7").Value = 9.5
$ws.Range("AX7").Value = 34

# Row 8
$ws.Range("H8").Value = 3.65
$ws.Range("I8").Value = 2.37
$ws.Range("K8").Value = 2.3
$ws.Range("L8").Value = 2.85
$ws.Range("Q8").Value = 1.55
$ws.Range("R8").Value = 2.15
$ws.Range("W8").Value = 11.5
$ws.Range("X8").Value = 15
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 15
$ws.Range("AD8").Value = 7.4
$ws.Range("AH8").Value = 11.25
$ws.Range("AL8").Value = 17
$ws.Range("AN8").Value = 4.75
$ws.Range("AO8").Value = 13
$ws.Range("AT8").Value = 2.95
$ws.Range("AU8").Value = 6.4
$ws.Range("AW8").Value = 4.6
$ws.Range("AX8").Value = 11.75
$ws.Range("AY8").Value = 16.5
$ws.Range("BA8").Value = 60
$ws.Range("BB8").Value = 150
